# Applies the cryptos-list refresh described by the commit:
#   "Updated cryptos list on Tue Aug  6 09:29:32 UTC 2024 with GitHub Actions"
#
# Column D holds price strings that look numeric (e.g. "479.32", "0.500",
# thousand-grouped "2.425.47"). Assigning them to Range.Value directly lets
# Excel auto-coerce plain decimals (e.g. "0.500" -> 0.5) and drop formatting,
# so we write them with a leading apostrophe (forces text, like typing it in
# the Excel UI) and then reset the cell Style back to "Normal" so no stray
# quote-prefix / number-format style sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''54.712.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.12%  '
# Row 3
$ws.Range("D3").Value = '''2.425.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.37%  '
# Row 4
$ws.Range("E4").Value = '  +0.05%  '
# Row 5
$ws.Range("D5").Value = '''479.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.11%  '
# Row 6
$ws.Range("D6").Value = '''137.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +13.35%  '
# Row 7
$ws.Range("E7").Value = '  -0.02%  '
# Row 8
$ws.Range("D8").Value = '''0.500'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.20%  '
# Row 9
$ws.Range("D9").Value = '''2.448.84'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.79%  '
# Row 10
$ws.Range("D10").Value = '''0.0967'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +11.67%  '
# Row 11
$ws.Range("D11").Value = '''5.45'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.69%  '
# Row 12
$ws.Range("E12").Value = '  +7.77%  '
# Row 13
$ws.Range("E13").Value = '  +1.44%  '
# Row 14
$ws.Range("D14").Value = '''2.862.22'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.52%  '
# Row 15
$ws.Range("D15").Value = '''54.851.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.38%  '
# Row 16
$ws.Range("D16").Value = '''20.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.97%  '
# Row 17
$ws.Range("D17").Value = '''0.0000134'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +14.09%  '
# Row 18
$ws.Range("D18").Value = '''2.451.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.66%  '
# Row 19
$ws.Range("D19").Value = '''4.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.88%  '
# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''313.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.44%  '
# Row 21
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '''9.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +11.19%  '
# Row 22
$ws.Range("E22").Value = '  -0.54%  '
# Row 23
$ws.Range("D23").Value = '''5.63'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.81%  '
# Row 24
$ws.Range("D24").Value = '''57.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.53%  '
# Row 25
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").Value = '''1.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.64%  '
# Row 26
$ws.Range("D26").Value = '''0.403'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.54%  '
# Row 27
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '''0.163'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.23%  '
# Row 28
$ws.Range("D28").Value = '''2.545.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.29%  '
# Row 29
$ws.Range("D29").Value = '''7.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.45%  '
# Row 30
$ws.Range("D30").Value = '''0.0₃0774'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +17.32%  '
# Row 31
$ws.Range("D31").Value = '''0.998'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.16%  '
# Row 32
$ws.Range("D32").Value = '''147.91'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.52%  '
# Row 33
$ws.Range("D33").Value = '''17.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.22%  '
# Row 34
$ws.Range("E34").Value = '  +9.73%  '
# Row 35
$ws.Range("D35").Value = '''5.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.81%  '
# Row 36
$ws.Range("E36").Value = '  +11.73%  '
# Row 37
$ws.Range("D37").Value = '''3.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.20%  '
# Row 38
$ws.Range("D38").Value = '''0.847'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.22%  '
# Row 39
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '''0.994'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.29%  '
# Row 40
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '''33.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.92%  '
# Row 41
$ws.Range("E41").Value = '  +9.58%  '
# Row 42
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").Value = '''0.0543'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.67%  '
# Row 43
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '''0.596'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.38%  '
# Row 44
$ws.Range("E44").Value = '  +10.54%  '
# Row 45
$ws.Range("D45").Value = '''10.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.07%  '
# Row 46
$ws.Range("D46").Value = '''255.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +26.52%  '
# Row 47
$ws.Range("D47").Value = '''4.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.02%  '
# Row 48
$ws.Range("D48").Value = '''0.0899'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.65%  '
# Row 49
$ws.Range("D49").Value = '''1.934.13'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.58%  '
# Row 50
$ws.Range("D50").Value = '''0.0222'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.34%  '
# Row 51
$ws.Range("D51").Value = '''17.04'
$ws.Range("D51").Style = "Normal"
